$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Data edit: column D held the literal placeholder "PEPTIDE" in every
#    data row; it now holds an actual RNA sequence, and the header
#    ("Sequence") becomes "RNA". Write the data rows first, then the
#    header, matching the order the strings were authored in (so the
#    shared-string table comes out in the same order as the target).
# ------------------------------------------------------------------
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 4).Value = "CCUGAACCCACUAUUGAUGAG"
}
$ws.Range("D1").Value = "RNA"

# ------------------------------------------------------------------
# 2) New formatting for the RNA column: a dedicated font (Arial Unicode
#    MS, 10pt) with vertically-centered text. Build it once on a scratch
#    cell, then copy/paste-special the formatting onto D2:D13 so the
#    whole block collapses onto a single new style record instead of one
#    per incremental property write.
# ------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.Font.Name = "Arial Unicode MS"
$scratch.Font.Size = 10
$scratch.VerticalAlignment = -4108
$scratch.Copy() | Out-Null
$ws.Range("D2:D13").PasteSpecial(-4122) | Out-Null
$scratch.Clear() | Out-Null

# ------------------------------------------------------------------
# 3) Column D needs to widen to fit the new text.
# ------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 28.33

# ------------------------------------------------------------------
# 4) Selection moves from H17 to D1.
# ------------------------------------------------------------------
$ws.Range("D1").Select() | Out-Null

# ------------------------------------------------------------------
# 5) Page setup: A4, portrait.
# ------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "edit applied"
